$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Letter date: "April 12, 2023" -> "April 26, 2023"
# ---------------------------------------------------------------------
$d.Content.Find.Execute("April 12, 2023", $true, $false, $false, $false, $false,
                         $true, 1, $false, "April 26, 2023", 2)

# ---------------------------------------------------------------------
# 2. Job title: "C++ Software Engineer" -> "Software Engineer (Data Platform)"
# ---------------------------------------------------------------------
$d.Content.Find.Execute("C++ Software Engineer", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Software Engineer (Data Platform)", 2)

# ---------------------------------------------------------------------
# 3. Company name: "iboss" -> "DoorDash", with Helvetica / black-colored run
#    (mirrors the pasted-text formatting seen in the target diff)
# ---------------------------------------------------------------------
$d.Content.Find.Execute("iboss", $true, $false, $false, $false, $false,
                         $true, 1, $false, "DoorDash", 2)

$companyRange = $d.Content
$companyRange.Find.MatchCase = $true
$companyRange.Find.Execute("DoorDash")
$companyRange.Font.Name = "Helvetica"
$companyRange.Font.Color = 0

# ---------------------------------------------------------------------
# 4. Skills paragraph: mention cryptography / OpenSSL before the MITRE
#    ATT&CK framework sentence.
# ---------------------------------------------------------------------
$d.Content.Find.Execute("I have worked extensively with the MITRE ATT&CK framework",
                         $true, $false, $false, $false, $false,
                         $true, 1, $false,
                         "I have worked extensively with cryptography, OpenSSL, MITRE ATT&CK framework",
                         2)
